$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B. This shifts the existing weekly
# columns (B:E -- both the header labels in row 1 and the per-analyst
# "UN"/rating data in rows 2-27) two weeks to the right, from B:E to E:H,
# freeing up B:D for the two newest snapshot weeks.
$ws.Range("B1:D1").EntireColumn.Insert()

# New header row: most-recent week first.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# The newly freed B:D columns (rows 2-27) need the same "UN" placeholder
# the other weekly columns use for analysts with no rating that week.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Add two new analyst-group rows at the bottom.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
